# Generate Report for Handback
# The handback transform failed because the handback file name did not
# match the expected handoff file name. Reflect this in the localization
# status report: update the Status for the affected file on both the
# zh-cn and de-de language sheets (and the roll-up Overview sheet), and
# record the error detail in the "Error Detail" column.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$zhError = "Handback file name: azgxcnol.d4f is different with handoff file name: 491be64c-002e-417b-9956-2b666f7340c5.d87886121f7405cdab8f6720cd72083a54263f29.zh-cn."
$deError = "Handback file name: azgxcnol.d4f is different with handoff file name: 491be64c-002e-417b-9956-2b666f7340c5.d87886121f7405cdab8f6720cd72083a54263f29.de-de."

# zh-cn sheet: row 3 corresponds to 491be64c-002e-417b-9956-2b666f7340c5
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("K3").Value = $zhError

# de-de sheet: row 3 corresponds to 491be64c-002e-417b-9956-2b666f7340c5
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("K3").Value = $deError

# Overview sheet rolls up the per-language status in columns B (zh-cn)
# and C (de-de) for the same file row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText
